$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 0.7441320593187497
$ws.Cells.Item(2, 3).Value = 0.9995917744288552
$ws.Cells.Item(2, 4).Value = 0.857921434222939
$ws.Cells.Item(2, 5).Value = 0.8166022341605902
$ws.Cells.Item(2, 6).Value = 0.7569633676092544
$ws.Cells.Item(2, 7).Value = 0.9277736891103059
$ws.Cells.Item(2, 8).Value = 0.9008891064586511
$ws.Cells.Item(2, 9).Value = 0.8280243026161923
$ws.Cells.Item(2, 10).Value = 0.9143653571891611
$ws.Cells.Item(2, 11).Value = 0.889498276776246
$ws.Cells.Item(2, 12).Value = 0.9017421379031961
$ws.Cells.Item(2, 13).Value = 0.9175424284798955

$ws.Cells.Item(3, 2).Value = 0.7371857750608142
$ws.Cells.Item(3, 3).Value = 0.999636073553902
$ws.Cells.Item(3, 4).Value = 0.8485058275767626
$ws.Cells.Item(3, 5).Value = 0.807192422457763
$ws.Cells.Item(3, 6).Value = 0.7990014460154241
$ws.Cells.Item(3, 7).Value = 0.9274308700684379
$ws.Cells.Item(3, 8).Value = 0.9002636175239414
$ws.Cells.Item(3, 9).Value = 0.8174478793824608
$ws.Cells.Item(3, 10).Value = 0.9143765340200378
$ws.Cells.Item(3, 11).Value = 0.8967602829358544
$ws.Cells.Item(3, 12).Value = 0.9027237532567611
$ws.Cells.Item(3, 13).Value = 0.9214871096761912

$ws.Cells.Item(4, 2).Value = 0.7398099268915896
$ws.Cells.Item(4, 3).Value = 0.9995858945192163
$ws.Cells.Item(4, 4).Value = 0.8514400619275129
$ws.Cells.Item(4, 5).Value = 0.8359218921183218
$ws.Cells.Item(4, 6).Value = 0.797036471722365
$ws.Cells.Item(4, 7).Value = 0.9218801669973544
$ws.Cells.Item(4, 8).Value = 0.9087481143028904
$ws.Cells.Item(4, 9).Value = 0.8182579025300972
$ws.Cells.Item(4, 10).Value = 0.9122665346310378
$ws.Cells.Item(4, 11).Value = 0.8890103590848583
$ws.Cells.Item(4, 12).Value = 0.8836648196396462
$ws.Cells.Item(4, 13).Value = 0.9283160395528551

$ws.Cells.Item(5, 2).Value = 0.730267954694243
$ws.Cells.Item(5, 3).Value = 0.9996031917291252
$ws.Cells.Item(5, 4).Value = 0.8383204064088128
$ws.Cells.Item(5, 5).Value = 0.8225630497759552
$ws.Cells.Item(5, 6).Value = 0.8004305912596401
$ws.Cells.Item(5, 7).Value = 0.9237238864446745
$ws.Cells.Item(5, 8).Value = 0.9083146157217512
$ws.Cells.Item(5, 9).Value = 0.823188185769152
$ws.Cells.Item(5, 10).Value = 0.9084157438329973
$ws.Cells.Item(5, 11).Value = 0.8903186412761445
$ws.Cells.Item(5, 12).Value = 0.8816213166489302
$ws.Cells.Item(5, 13).Value = 0.9146501946149194

$ws.Cells.Item(6, 2).Value = 0.7320967597143977
$ws.Cells.Item(6, 3).Value = 0.9995634309809356
$ws.Cells.Item(6, 4).Value = 0.8536349129667387
$ws.Cells.Item(6, 5).Value = 0.8165106129227874
$ws.Cells.Item(6, 6).Value = 0.781844473007712
$ws.Cells.Item(6, 7).Value = 0.9253016178542675
$ws.Cells.Item(6, 8).Value = 0.902256929136013
$ws.Cells.Item(6, 9).Value = 0.8344343434440493
$ws.Cells.Item(6, 10).Value = 0.9124027056872186
$ws.Cells.Item(6, 11).Value = 0.8976092033121996
$ws.Cells.Item(6, 12).Value = 0.8635933085024403
$ws.Cells.Item(6, 13).Value = 0.9162983367179868

$ws.Cells.Item(7, 2).Value = 0.7349064742215674
$ws.Cells.Item(7, 3).Value = 0.9996595075627055
$ws.Cells.Item(7, 4).Value = 0.8560133022427803
$ws.Cells.Item(7, 5).Value = 0.7940797151355956
$ws.Cells.Item(7, 6).Value = 0.7870629820051414
$ws.Cells.Item(7, 7).Value = 0.9181520660800552
$ws.Cells.Item(7, 8).Value = 0.8903934890076837
$ws.Cells.Item(7, 9).Value = 0.8164504843323459
$ws.Cells.Item(7, 10).Value = 0.9114347921332993
$ws.Cells.Item(7, 11).Value = 0.8886183327316622
$ws.Cells.Item(7, 12).Value = 0.9172702377894389
$ws.Cells.Item(7, 13).Value = 0.9177874139436655

$ws.Cells.Item(8, 2).Value = 0.7342452361298741
$ws.Cells.Item(8, 3).Value = 0.9995815559451139
$ws.Cells.Item(8, 4).Value = 0.8422279639998783
$ws.Cells.Item(8, 5).Value = 0.7712062401827455
$ws.Cells.Item(8, 6).Value = 0.7904330012853471
$ws.Cells.Item(8, 7).Value = 0.9188601693696803
$ws.Cells.Item(8, 8).Value = 0.902832900196725
$ws.Cells.Item(8, 9).Value = 0.812349441871951
$ws.Cells.Item(8, 10).Value = 0.9139270391382812
$ws.Cells.Item(8, 11).Value = 0.8924825421357822
$ws.Cells.Item(8, 12).Value = 0.9840980606216285
$ws.Cells.Item(8, 13).Value = 0.9186996414332703

$ws.Cells.Item(9, 2).Value = 0.7358446878715359
$ws.Cells.Item(9, 3).Value = 0.9996131533236104
$ws.Cells.Item(9, 4).Value = 0.8533752967043
$ws.Cells.Item(9, 5).Value = 0.7925066250695895
$ws.Cells.Item(9, 6).Value = 0.7904105077120822
$ws.Cells.Item(9, 7).Value = 0.922302288216011
$ws.Cells.Item(9, 8).Value = 0.9000720615376863
$ws.Cells.Item(9, 9).Value = 0.8150206925842044
$ws.Cells.Item(9, 10).Value = 0.9155828866326593
$ws.Cells.Item(9, 11).Value = 0.8943629994810586
$ws.Cells.Item(9, 12).Value = 0.8998672067080108
$ws.Cells.Item(9, 13).Value = 0.9322795657848255

$ws.Cells.Item(10, 2).Value = 0.7365376741497507
$ws.Cells.Item(10, 3).Value = 0.9995925450966234
$ws.Cells.Item(10, 4).Value = 0.8497639679255045
$ws.Cells.Item(10, 5).Value = 0.8119039281441031
$ws.Cells.Item(10, 6).Value = 0.7905004820051413
$ws.Cells.Item(10, 7).Value = 0.9241075921612607
$ws.Cells.Item(10, 8).Value = 0.9072890744892167
$ws.Cells.Item(10, 9).Value = 0.8147184953363613
$ws.Cells.Item(10, 10).Value = 0.9079098059553137
$ws.Cells.Item(10, 11).Value = 0.8899812588840506
$ws.Cells.Item(10, 12).Value = 0.9361674800924736
$ws.Cells.Item(10, 13).Value = 0.9160664469569739

$ws.Cells.Item(11, 2).Value = 0.749696202697939
$ws.Cells.Item(11, 3).Value = 0.9995657429842402
$ws.Cells.Item(11, 4).Value = 0.8437148139131494
$ws.Cells.Item(11, 5).Value = 0.8132301067387421
$ws.Cells.Item(11, 6).Value = 0.80188544344473
$ws.Cells.Item(11, 7).Value = 0.91993400621118
$ws.Cells.Item(11, 8).Value = 0.9034705514162763
$ws.Cells.Item(11, 9).Value = 0.8160709361196337
$ws.Cells.Item(11, 10).Value = 0.9072628537280691
$ws.Cells.Item(11, 11).Value = 0.8896445815753256
$ws.Cells.Item(11, 12).Value = 0.900788044475432
$ws.Cells.Item(11, 13).Value = 0.917396139897879

$ws.Cells.Item(12, 2).Value = 0.711999062881115
$ws.Cells.Item(12, 3).Value = 0.9995603768531134
$ws.Cells.Item(12, 4).Value = 0.8530790678920299
$ws.Cells.Item(12, 5).Value = 0.7960006553247857
$ws.Cells.Item(12, 6).Value = 0.7591171272493573
$ws.Cells.Item(12, 7).Value = 0.9226976758396594
$ws.Cells.Item(12, 8).Value = 0.9060389653544285
$ws.Cells.Item(12, 9).Value = 0.822826894305912
$ws.Cells.Item(12, 10).Value = 0.9076808672028565
$ws.Cells.Item(12, 11).Value = 0.8883380620924618
$ws.Cells.Item(12, 12).Value = 0.9163952698983524
$ws.Cells.Item(12, 13).Value = 0.9175938530687704

$ws.Cells.Item(13, 2).Value = 0.7319434926070516
$ws.Cells.Item(13, 3).Value = 0.9996537132828186
$ws.Cells.Item(13, 4).Value = 0.8516430951583944
$ws.Cells.Item(13, 5).Value = 0.8135189465731717
$ws.Cells.Item(13, 6).Value = 0.7748481683804627
$ws.Cells.Item(13, 7).Value = 0.919758552953186
$ws.Cells.Item(13, 8).Value = 0.902291244153959
$ws.Cells.Item(13, 9).Value = 0.8046566241732815
$ws.Cells.Item(13, 10).Value = 0.9172728234612112
$ws.Cells.Item(13, 11).Value = 0.8915909642156088
$ws.Cells.Item(13, 12).Value = 0.8830272283585923
$ws.Cells.Item(13, 13).Value = 0.9264139685792567

$ws.Cells.Item(14, 2).Value = 0.7476632239954983
$ws.Cells.Item(14, 3).Value = 0.9996385568167105
$ws.Cells.Item(14, 4).Value = 0.842816617723137
$ws.Cells.Item(14, 5).Value = 0.7981576366689934
$ws.Cells.Item(14, 6).Value = 0.7925401670951158
$ws.Cells.Item(14, 7).Value = 0.9224292166292845
$ws.Cells.Item(14, 8).Value = 0.9003070542555182
$ws.Cells.Item(14, 9).Value = 0.8133871939472149
$ws.Cells.Item(14, 10).Value = 0.9117086244897776
$ws.Cells.Item(14, 11).Value = 0.8946735887051285
$ws.Cells.Item(14, 12).Value = 0.9193630508972148
$ws.Cells.Item(14, 13).Value = 0.9198859204586178

$ws.Cells.Item(15, 2).Value = 0.722258106187831
$ws.Cells.Item(15, 3).Value = 0.999596969300478
$ws.Cells.Item(15, 4).Value = 0.8443286684823954
$ws.Cells.Item(15, 5).Value = 0.8020973498945191
$ws.Cells.Item(15, 6).Value = 0.812538560411311
$ws.Cells.Item(15, 7).Value = 0.9209465130406027
$ws.Cells.Item(15, 8).Value = 0.909991273560626
$ws.Cells.Item(15, 9).Value = 0.8273204568799599
$ws.Cells.Item(15, 10).Value = 0.9126310855981318
$ws.Cells.Item(15, 11).Value = 0.8908037386340563
$ws.Cells.Item(15, 12).Value = 0.8778290246229495
$ws.Cells.Item(15, 13).Value = 0.9230739255614702

$ws.Cells.Item(16, 2).Value = 0.7306199216586128
$ws.Cells.Item(16, 3).Value = 0.9996075017599769
$ws.Cells.Item(16, 4).Value = 0.8425089772070236
$ws.Cells.Item(16, 5).Value = 0.8149491696242986
$ws.Cells.Item(16, 6).Value = 0.7965223329048843
$ws.Cells.Item(16, 7).Value = 0.9182062072262479
$ws.Cells.Item(16, 8).Value = 0.8973029698996482
$ws.Cells.Item(16, 9).Value = 0.8196406110822313
$ws.Cells.Item(16, 10).Value = 0.9169133020680118
$ws.Cells.Item(16, 11).Value = 0.8882880011732588
$ws.Cells.Item(16, 12).Value = 0.9014588914168287
$ws.Cells.Item(16, 13).Value = 0.9215698361887288

$ws.Cells.Item(17, 2).Value = 0.7289224884447549
$ws.Cells.Item(17, 3).Value = 0.9996666719186248
$ws.Cells.Item(17, 4).Value = 0.8550414054654454
$ws.Cells.Item(17, 5).Value = 0.8197554023531129
$ws.Cells.Item(17, 6).Value = 0.799146850899743
$ws.Cells.Item(17, 7).Value = 0.9170142034017714
$ws.Cells.Item(17, 8).Value = 0.8944860978568751
$ws.Cells.Item(17, 9).Value = 0.8222777504994181
$ws.Cells.Item(17, 10).Value = 0.9130971594456888
$ws.Cells.Item(17, 11).Value = 0.8939878951287201
$ws.Cells.Item(17, 12).Value = 0.9196508843712158
$ws.Cells.Item(17, 13).Value = 0.9264797665004258

$ws.Cells.Item(18, 2).Value = 0.7318367530144356
$ws.Cells.Item(18, 3).Value = 0.9995940864321599
$ws.Cells.Item(18, 4).Value = 0.8578776893582059
$ws.Cells.Item(18, 5).Value = 0.8392140450698844
$ws.Cells.Item(18, 6).Value = 0.7947188303341902
$ws.Cells.Item(18, 7).Value = 0.9273749317057741
$ws.Cells.Item(18, 8).Value = 0.9025531676453684
$ws.Cells.Item(18, 9).Value = 0.8131753195651819
$ws.Cells.Item(18, 10).Value = 0.9112190792973798
$ws.Cells.Item(18, 11).Value = 0.8927416602739108
$ws.Cells.Item(18, 12).Value = 0.8818116766357197
$ws.Cells.Item(18, 13).Value = 0.9197124822489346

$ws.Cells.Item(19, 2).Value = 0.7467217260503722
$ws.Cells.Item(19, 3).Value = 0.9995630599186769
$ws.Cells.Item(19, 4).Value = 0.8568192538267247
$ws.Cells.Item(19, 5).Value = 0.8210559580592137
$ws.Cells.Item(19, 6).Value = 0.8238174807197942
$ws.Cells.Item(19, 7).Value = 0.9207966701173222
$ws.Cells.Item(19, 8).Value = 0.8966496814567289
$ws.Cells.Item(19, 9).Value = 0.815009642446366
$ws.Cells.Item(19, 10).Value = 0.9098640748340987
$ws.Cells.Item(19, 11).Value = 0.8995161012838157
$ws.Cells.Item(19, 12).Value = 0.9012008733624455
$ws.Cells.Item(19, 13).Value = 0.9254161399106554

$ws.Cells.Item(20, 2).Value = 0.7368214919896041
$ws.Cells.Item(20, 3).Value = 0.9995937724564025
$ws.Cells.Item(20, 4).Value = 0.8428689213657525
$ws.Cells.Item(20, 5).Value = 0.8125134423108428
$ws.Cells.Item(20, 6).Value = 0.7949558161953727
$ws.Cells.Item(20, 7).Value = 0.9250841546612606
$ws.Cells.Item(20, 8).Value = 0.902677396697679
$ws.Cells.Item(20, 9).Value = 0.8041209327085042
$ws.Cells.Item(20, 10).Value = 0.9113880357241321
$ws.Cells.Item(20, 11).Value = 0.8884343059723382
$ws.Cells.Item(20, 12).Value = 0.8820559337272027
$ws.Cells.Item(20, 13).Value = 0.9216573538120314

$ws.Cells.Item(21, 2).Value = 0.7334247097230463
$ws.Cells.Item(21, 3).Value = 0.9995777596927741
$ws.Cells.Item(21, 4).Value = 0.8469409976872281
$ws.Cells.Item(21, 5).Value = 0.8061279079066799
$ws.Cells.Item(21, 6).Value = 0.7292424485861182
$ws.Cells.Item(21, 7).Value = 0.9250412461180124
$ws.Cells.Item(21, 8).Value = 0.90579745712686
$ws.Cells.Item(21, 9).Value = 0.8140425151651132
$ws.Cells.Item(21, 10).Value = 0.9148016261543803
$ws.Cells.Item(21, 11).Value = 0.8970500722005369
$ws.Cells.Item(21, 12).Value = 0.8840891159957432
$ws.Cells.Item(21, 13).Value = 0.9177059650606652

